$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.862.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.73%  "

# Row 3: update D3, E3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.374.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.97%  "

# Row 4: update E4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.83%  "

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.20%  "

# Row 7: update E7
$ws.Range("E7").Value = "  -2.49%  "

# Row 8: update E8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9: update D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.57%  "

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.16%  "

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.11%  "

# Row 12: update D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.45%  "

# Row 13: update E13
$ws.Range("E13").Value = "  -3.54%  "

# Row 14: update E14
$ws.Range("E14").Value = "  -0.25%  "

# Row 15: update D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.08%  "

# Row 16: update D16, E16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.732.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.09%  "

# Row 17: update D17, E17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.349.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.85%  "

# Row 18: update D18, E18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.834.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.49%  "

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000107"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.00%  "

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.85%  "

# Row 23: update D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.04%  "

# Row 24: update E24
$ws.Range("E24").Value = "  -3.93%  "

# Row 25: update D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.29%  "

# Row 26: update E26
$ws.Range("E26").Value = "  +0.17%  "

# Row 27: update D27, E27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.47%  "

# Row 28: update D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "

# Row 29: update E29
$ws.Range("E29").Value = "  +2.87%  "

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.11%  "

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0898"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.21%  "

# Row 33: update D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.31%  "

# Row 34: update E34
$ws.Range("E34").Value = "  -7.21%  "

# Row 35: update E35
$ws.Range("E35").Value = "  -2.58%  "

# Row 36: update D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.119"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.22%  "

# Row 37: update D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.07%  "

# Row 38: update D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0366"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.99%  "

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.12%  "

# Row 41: update B41, C41, D41, E41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.50%  "

# Row 42: update B42, C42, D42, E42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.242"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.05%  "

# Row 43: update D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.87%  "

# Row 44: update E44
$ws.Range("E44").Value = "  -0.14%  "

# Row 45: update B45, C45, D45, E45
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "

# Row 46: update B46, C46, D46, E46
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.87%  "

# Row 47: update D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "113.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.14%  "

# Row 48: update D48, E48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.28%  "

# Row 49: update D49, E49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "

# Row 50: update D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.40%  "

# Row 51: update D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.30%  "
